# fix(module3): use uncon_planned_qty for future production; keep produced for today
# Populates DeliveryPlan / VehicleLog / TruckUsageLog / ValidationLog sheets
# with the simulation output rows, re-labels headers, and formats the new
# date columns.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("DeliveryPlan")
$ws2 = $wb.Worksheets.Item("VehicleLog")
$ws3 = $wb.Worksheets.Item("TruckUsageLog")

# ---------------------------------------------------------------------
# DeliveryPlan (sheet1)
# ---------------------------------------------------------------------

$ws1.Range("A1").Value = "vehicle_uid"
$ws1.Range("B1").Value = "ori_deployment_uid"
$ws1.Range("C1").Value = "material"
$ws1.Range("D1").Value = "sending"
$ws1.Range("E1").Value = "receiving"
$ws1.Range("F1").Value = "planned_deployment_date"
$ws1.Range("G1").Value = "actual_ship_date"
$ws1.Range("H1").Value = "actual_delivery_date"
$ws1.Range("I1").Value = "delivery_qty"
$ws1.Range("J1").Value = "truck_type"
$ws1.Range("K1").Value = "truck_load_pct"
$ws1.Range("L1").Value = "WFR"
$ws1.Range("M1").Value = "VFR"

# VehicleLog!A1 already carries the bold/centered/bordered header style from
# the original workbook - clone it onto the new DeliveryPlan header row
# instead of re-deriving it property-by-property (which mints spurious
# intermediate cellXfs entries).
$ws2.Range("A1").Copy()
$ws1.Range("A1:M1").PasteSpecial(-4122)

$ws1.Range("A2").Value = "20240103-PLANT_001-DC_001-LARGE-#1"
$ws1.Range("B2").Value = "MAT_B|PLANT_001|DC_001|2024-01-06|net demand for forecast|000046"
$ws1.Range("C2").Value = "MAT_B"
$ws1.Range("D2").Value = "PLANT_001"
$ws1.Range("E2").Value = "DC_001"
$ws1.Range("F2").Value = 45297
$ws1.Range("G2").Value = 45294
$ws1.Range("H2").Value = 45296
$ws1.Range("I2").Value = 35
$ws1.Range("J2").Value = "LARGE"
$ws1.Range("K2").Value = 0.76
$ws1.Range("L2").Value = 0.76
$ws1.Range("M2").Value = 0.7125

$ws1.Range("A3").Value = "20240103-PLANT_001-DC_001-LARGE-#1"
$ws1.Range("B3").Value = "MAT_B|PLANT_001|DC_001|2024-01-01|net demand for safety|000023"
$ws1.Range("C3").Value = "MAT_B"
$ws1.Range("D3").Value = "PLANT_001"
$ws1.Range("E3").Value = "DC_001"
$ws1.Range("F3").Value = 45292
$ws1.Range("G3").Value = 45294
$ws1.Range("H3").Value = 45296
$ws1.Range("I3").Value = 60
$ws1.Range("J3").Value = "LARGE"
$ws1.Range("K3").Value = 0.76
$ws1.Range("L3").Value = 0.76
$ws1.Range("M3").Value = 0.7125

$ws1.Range("F2:H3").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# VehicleLog (sheet2)
# ---------------------------------------------------------------------

$ws2.Range("A1").Value = "date"
$ws2.Range("B1").Value = "sending"
$ws2.Range("C1").Value = "receiving"
$ws2.Range("D1").Value = "truck_type"
$ws2.Range("E1").Value = "vehicle_no"
$ws2.Range("F1").Value = "vehicle_uid"
$ws2.Range("G1").Value = "total_units"
$ws2.Range("H1").Value = "total_weight"
$ws2.Range("I1").Value = "total_volume"
$ws2.Range("J1").Value = "WFR"
$ws2.Range("K1").Value = "VFR"
$ws2.Range("L1").Value = "trigger"

$ws2.Range("A2").Value = 45294
$ws2.Range("B2").Value = "PLANT_001"
$ws2.Range("C2").Value = "DC_001"
$ws2.Range("D2").Value = "LARGE"
$ws2.Range("E2").Value = 1
$ws2.Range("F2").Value = "20240103-PLANT_001-DC_001-LARGE-#1"
$ws2.Range("G2").Value = 95
$ws2.Range("H2").Value = 76
$ws2.Range("I2").Value = 142.5
$ws2.Range("J2").Value = 0.76
$ws2.Range("K2").Value = 0.7125
$ws2.Range("L2").Value = "threshold"

$ws2.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# TruckUsageLog (sheet3)
# ---------------------------------------------------------------------

$ws3.Range("A1").Value = "date"
$ws3.Range("B1").Value = "sending"
$ws3.Range("C1").Value = "receiving"
$ws3.Range("D1").Value = "truck_type"
$ws3.Range("E1").Value = "truck_used"

$ws3.Range("A2").Value = 45294
$ws3.Range("B2").Value = "PLANT_001"
$ws3.Range("C2").Value = "DC_001"
$ws3.Range("D2").Value = "LARGE"
$ws3.Range("E2").Value = 1

$ws3.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

